# Insert a new data row at row 56 (pushing the existing rows 56-80 down to
# 57-81) and populate it with the new "Haba" price record for
# Macroferia Regional de Talca / Provincia del Elquí dated 44806.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(56).Insert()

$ws.Cells.Item(56, 1).Value = 5
$ws.Cells.Item(56, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(56, 3).Value = "Maule"
$ws.Cells.Item(56, 4).Value = 44806
$ws.Cells.Item(56, 5).Value = 7
$ws.Cells.Item(56, 6).Value = 100112026
$ws.Cells.Item(56, 7).Value = "Haba"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 200
$ws.Cells.Item(56, 11).Value = 12000
$ws.Cells.Item(56, 12).Value = 12000
$ws.Cells.Item(56, 13).Value = 12000
$ws.Cells.Item(56, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(56, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(56, 16).Value = 480
$ws.Cells.Item(56, 17).Value = 25
$ws.Cells.Item(56, 18).Value = "Hortaliza"
